# Generate Report for Handback
# Updates the handback-status report with the latest handoff/handback
# file names and timestamps:
#   647be52a-93d1-4692-8474-8d5d0c3bc3e3 -> 90ae4d30-adb9-4c90-a532-eb17ba8b0997
#   a96955cb-9400-4f1d-ac1c-668f24ff21f5 -> ffff655aa22f-b345-4971-a2a7-e9228d8c6fe2

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview.Range("A2").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.md"
$wsOverview.Range("B2").Value = "e2e\90ae4d30-adb9-4c90-a532-eb17ba8b0997.md"
$wsOverview.Range("G2").Value = "2016-09-03 11:06:51"

$wsOverview.Range("A3").Value = "ffff655aa22f-b345-4971-a2a7-e9228d8c6fe2.md"
$wsOverview.Range("B3").Value = "e2e\ffff655aa22f-b345-4971-a2a7-e9228d8c6fe2.md"
$wsOverview.Range("G3").Value = "2016-09-03 11:06:51"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZhCn.Range("A2").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.md"
$wsZhCn.Range("G2").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.bdc58107db3ad851ca5abdf44805a75182a73397.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-03 11:06:47"
$wsZhCn.Range("I2").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.md"
$wsZhCn.Range("J2").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.bdc58107db3ad851ca5abdf44805a75182a73397.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-03 11:07:26"

$wsZhCn.Range("A3").Value = "ffff655aa22f-b345-4971-a2a7-e9228d8c6fe2.md"
$wsZhCn.Range("G3").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.bdc58107db3ad851ca5abdf44805a75182a73397.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-03 11:06:47"
$wsZhCn.Range("I3").Value = "ffff655aa22f-b345-4971-a2a7-e9228d8c6fe2.md"
$wsZhCn.Range("J3").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.bdc58107db3ad851ca5abdf44805a75182a73397.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-03 11:07:26"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDeDe.Range("A2").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.md"
$wsDeDe.Range("G2").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.bdc58107db3ad851ca5abdf44805a75182a73397.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-03 11:06:51"
$wsDeDe.Range("I2").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.md"
$wsDeDe.Range("J2").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.bdc58107db3ad851ca5abdf44805a75182a73397.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-03 11:07:33"

$wsDeDe.Range("A3").Value = "ffff655aa22f-b345-4971-a2a7-e9228d8c6fe2.md"
$wsDeDe.Range("G3").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.bdc58107db3ad851ca5abdf44805a75182a73397.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-03 11:06:51"
$wsDeDe.Range("I3").Value = "ffff655aa22f-b345-4971-a2a7-e9228d8c6fe2.md"
$wsDeDe.Range("J3").Value = "90ae4d30-adb9-4c90-a532-eb17ba8b0997.bdc58107db3ad851ca5abdf44805a75182a73397.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-03 11:07:33"
